$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style from an existing header cell (e.g. H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data rows for columns I (I0) and J (IF)
$values = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(8, 8)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(6, 7)
    8  = @(8, 8)
    9  = @(6, 6)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(3, 3)
    13 = @(4, 4)
    14 = @(3, 3)
    15 = @(5, 5)
    16 = @(5, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
